$wb = $excel.ActiveWorkbook

# Rename the "Art._PTA" sheet to "Art_PTA"
$ws = $wb.Worksheets.Item("Art._PTA")
$ws.Name = "Art_PTA"

# Make Art_PTA the active sheet and set its selection to H14
$ws.Activate()
$ws.Range("H14").Select()
